$d = $word.ActiveDocument

# Shared rPr block used by every body run in this document.
$rPr = '<w:rPr><w:rFonts w:asciiTheme="majorBidi" w:hAnsiTheme="majorBidi" w:cstheme="majorBidi"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>'

$pkgOpen = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$pkgClose = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# ---------------------------------------------------------------------
# 1) "Vision" paragraph: split the run containing "Unity" so the word is
#    wrapped in spellcheck proofErr markers.
# ---------------------------------------------------------------------
$find1 = $d.Content
$find1.Find.Text = "Nuestra visi" + [char]0x00F3 + "n es poder crear un juego al estilo de " + [char]0x201C + "Elige tu propia aventura" + [char]0x201D + ", con varios finales, desarroll" + [char]0x00E1 + "ndolo en Unity y para poder ambientarlo usar VR."
$find1.Find.Forward = $true
$find1.Find.Wrap = 0
$find1.Find.MatchCase = $true
if (-not $find1.Find.Execute()) { throw "paragraph 1 text not found" }

$target1 = $d.Range($find1.Start, $find1.End)

$body1 = "<w:r>$rPr<w:t xml:space=`"preserve`">Nuestra visi" + [char]0x00F3 + "n es poder crear un juego al estilo de " + [char]0x201C + "Elige tu propia aventura" + [char]0x201D + ", con varios finales, desarroll" + [char]0x00E1 + "ndolo en </w:t></w:r>" +
         '<w:proofErr w:type="spellStart"/>' +
         "<w:r>$rPr<w:t>Unity</w:t></w:r>" +
         '<w:proofErr w:type="spellEnd"/>' +
         "<w:r>$rPr<w:t xml:space=`"preserve`"> y para poder ambientarlo usar VR.</w:t></w:r>"

$xml1 = $pkgOpen + "<w:p>" + $body1 + "</w:p>" + $pkgClose
$target1.InsertXML($xml1)

# ---------------------------------------------------------------------
# 2) "El mercado ahora" paragraph: split the run containing "stores" so
#    the word is wrapped in spellcheck proofErr markers, with the
#    _GoBack bookmark re-homed inside it (between "store" and "s").
#    Remove the pre-existing _GoBack bookmark first so the name stays
#    unique while we recreate it in its new home.
# ---------------------------------------------------------------------
$oldBookmark = $d.Bookmarks("_GoBack")
$oldBookmark.Delete()

$find2 = $d.Content
$find2.Find.Text = " 10 juegos de este g" + [char]0x00E9 + "nero, pero ninguno de estos es compatible con VR. La mayor" + [char]0x00ED + "a tiene una gr" + [char]0x00E1 + "fica 3D, son en formato de primera (FP) o tercera persona (TP). Casi todos estos se venden en una plataforma virtual llamada " + [char]0x201C + "STEAM" + [char]0x201D + ", algunos est" + [char]0x00E1 + "n disponibles en los stores de consolas como PlayStation, etc. Por " + [char]0x00FA + "ltimo, los precios de estos rondan entre los 15 a 60 USD."
$find2.Find.Forward = $true
$find2.Find.Wrap = 0
$find2.Find.MatchCase = $true
if (-not $find2.Find.Execute()) { throw "paragraph 2 text not found" }

$target2 = $d.Range($find2.Start, $find2.End)

$body2 = "<w:r>$rPr<w:t xml:space=`"preserve`"> 10 juegos de este g" + [char]0x00E9 + "nero, pero ninguno de estos es compatible con VR. La mayor" + [char]0x00ED + "a tiene una gr" + [char]0x00E1 + "fica 3D, son en formato de primera (FP) o tercera persona (TP). Casi todos estos se venden en una plataforma virtual llamada " + [char]0x201C + "STEAM" + [char]0x201D + ", algunos est" + [char]0x00E1 + "n disponibles en los </w:t></w:r>" +
         '<w:proofErr w:type="spellStart"/>' +
         "<w:r>$rPr<w:t>store</w:t></w:r>" +
         '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
         "<w:r>$rPr<w:t>s</w:t></w:r>" +
         '<w:proofErr w:type="spellEnd"/>' +
         "<w:r>$rPr<w:t xml:space=`"preserve`"> de consolas como PlayStation, etc. Por " + [char]0x00FA + "ltimo, los precios de estos rondan entre los 15 a 60 USD.</w:t></w:r>"

$xml2 = $pkgOpen + "<w:p>" + $body2 + "</w:p>" + $pkgClose
$target2.InsertXML($xml2)

# ---------------------------------------------------------------------
# 3) "Queremos" paragraph: finish the sentence with a new run. Rebuild
#    the whole paragraph (minus its end-of-paragraph mark) so the now-
#    empty old bookmark location is cleanly dropped along with it.
# ---------------------------------------------------------------------
$find3 = $d.Content
$find3.Find.Text = "Queremos que este juego "
$find3.Find.Forward = $true
$find3.Find.Wrap = 0
$find3.Find.MatchCase = $true
if (-not $find3.Find.Execute()) { throw "paragraph 3 text not found" }

$target3 = $d.Range($find3.Start, $find3.End)

$body3 = "<w:r w:rsidRPr=`"004C1792`">$rPr<w:t>Queremos</w:t></w:r>" +
         "<w:r>$rPr<w:t xml:space=`"preserve`"> </w:t></w:r>" +
         "<w:r w:rsidR=`"00A80C13`">$rPr<w:t xml:space=`"preserve`">que este juego </w:t></w:r>" +
         "<w:r>$rPr<w:t>sea para pc, siendo complementado con VR.</w:t></w:r>"

$xml3 = $pkgOpen + "<w:p>" + $body3 + "</w:p>" + $pkgClose
$target3.InsertXML($xml3)
